$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 (bec9db9b...) currently shares the same handoff/handback
# datetimes as row 2. Give it its own distinct timestamps.
$wsZh.Range("E3").Value = "2016-03-11 08:28:42"
$wsZh.Range("H3").Value = "2016-03-11 08:28:59"

# de-de: same fix for row 3 (bec9db9b...).
$wsDe.Range("E3").Value = "2016-03-11 08:28:45"
$wsDe.Range("H3").Value = "2016-03-11 08:29:05"
